# Add a new kinetic model row (Streptococcus pneumoniae / NCTC7465 already had
# its name/strain filled in, but was missing the model id, year and page
# counts). The model id "spnLHP26" is a brand-new value not previously in the
# shared-string table, so just assigning it to A20 will append it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "spnLHP26"
$ws.Range("B20").Value = 2026
$ws.Range("E20").Value = 836
$ws.Range("F20").Value = 460
$ws.Range("G20").Value = 1014

# Move/record the active selection on the sheet, as captured in the saved file.
$ws.Range("A21").Select()
